$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.333.89'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.873.08'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7121'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.64'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3110'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.84%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07783'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.08'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +1.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08410'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.879.51'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.235'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7115'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.06'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.343.95'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.096'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +2.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008232'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +5.32%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '239.92'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.09%  '
$ws.Range("E20").Value = '  +0.94%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.121.54'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.758'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1592'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.20%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.90'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.039'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.48'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.509'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.86%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.414'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("E31").Value = '  -1.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.317'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +2.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05292'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +3.10%  '
$ws.Range("E34").Value = '  +1.15%  '
$ws.Range("E35").Value = '  +1.36%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7391'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -7.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.701'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +0.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01873'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +1.78%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.217.54'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +4.41%  '
$ws.Range("E40").Value = '  +1.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.556'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +6.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '110.94'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +8.65%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '72.81'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.9999'
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.021.08'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("E47").Value = '  +1.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5193'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.09%  '
$ws.Range("E49").Value = '  +1.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.397'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("E51").Value = '  +1.07%  '
